$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "91.489.49"
$ws.Range("E2").Value = "  -3.56%  "
$ws.Range("D3").Value = "3.328.77"
$ws.Range("E3").Value = "  -4.82%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "615.17"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.80%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.42"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.387"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.32%  "
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.960"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.20%  "
$ws.Range("D11").Value = "3.330.35"
$ws.Range("E11").Value = "  -4.81%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.26"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.26%  "
$ws.Range("E13").Value = "  -2.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.97"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.03%  "
$ws.Range("D15").Value = "91.352.61"
$ws.Range("E15").Value = "  -3.52%  "
$ws.Range("D16").Value = "3.948.80"
$ws.Range("E16").Value = "  -4.57%  "
$ws.Range("E17").Value = "  -5.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.09"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.86%  "
$ws.Range("D19").Value = "3.323.30"
$ws.Range("E19").Value = "  -4.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.43"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "490.68"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.450"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -11.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000182"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.11"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -9.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "89.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.80"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.96%  "
$ws.Range("D29").Value = "3.503.13"
$ws.Range("E29").Value = "  -4.67%  "
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("E31").Value = "  -7.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.138"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.60"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.997"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("E35").Value = "  -7.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "28.19"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.525"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -9.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "557.69"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.34%  "
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.34"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.37"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.90%  "
$ws.Range("E42").Value = "  -2.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.866"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "23.71"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.56%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.65"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.36%  "
$ws.Range("B46").Value = "ImmutableX"
$ws.Range("C46").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.11%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0412"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.11"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.27%  "
$ws.Range("B50").Value = "OKB"
$ws.Range("C50").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "51.76"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.32%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.94"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.17%  "

Write-Host "Updated 98 cells"
